$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "15-03-2025 00:00:00"
$ws.Range("K1").Value = "15-03-2025 00:00:00"
$ws.Range("M1").Value = 45731
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 5233.14
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 667.84
$ws.Range("F22").Value = 81
$ws.Range("G22").Value = 11792.79
$ws.Range("F23").Value = 69
$ws.Range("G23").Value = 10881.3
$ws.Range("B25").Value = 81732.95
$ws.Range("F59").Value = 32
$ws.Range("G59").Value = 4012.8
$ws.Range("F69").Value = 314
$ws.Range("G69").Value = 35224.52
$ws.Range("F71").Value = 216
$ws.Range("G71").Value = 9698.4
$ws.Range("F73").Value = 8
$ws.Range("G73").Value = 830.96
$ws.Range("F76").Value = 8
$ws.Range("G76").Value = 2100.64
$ws.Range("F82").Value = 48
$ws.Range("G82").Value = 793.92
$ws.Range("F89").Value = 10
$ws.Range("G89").Value = 505.5
$ws.Range("F92").Value = 93
$ws.Range("G92").Value = 4157.1
$ws.Range("B95").Value = 130227.99
$ws.Range("B99").Value = 48264
$ws.Range("C99").Value = "BLUE-Inverter 1.5 ton 5 star Split AC"
$ws.Range("D99").Value = 32287.23
$ws.Range("E99").Value = 41844.24
$ws.Range("F99").Value = 1
$ws.Range("G99").Value = 32287.23
$ws.Range("B100").Value = 54863
$ws.Range("C100").Value = "BLUE-Inverter 1.5 Ton 5 Star Split Ac"
$ws.Range("D100").Value = 32143.58
$ws.Range("E100").Value = 41658.07
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("F122").Value = 39
$ws.Range("G122").Value = 2399.28
$ws.Range("B124").Value = 4571.99
$ws.Range("F144").Value = 15
$ws.Range("G144").Value = 1509.45
$ws.Range("F145").Value = 12
$ws.Range("G145").Value = 612.24
$ws.Range("F146").Value = 90
$ws.Range("G146").Value = 6791.4
$ws.Range("F151").Value = 2
$ws.Range("G151").Value = 66.02
$ws.Range("F154").Value = 80
$ws.Range("G154").Value = 4059.2
$ws.Range("F159").Value = 2
$ws.Range("G159").Value = 108.82
$ws.Range("B168").Value = 32481.54
$ws.Range("F173").Value = 16
$ws.Range("G173").Value = 1225.6
$ws.Range("F174").Value = 38
$ws.Range("G174").Value = 4251.44
$ws.Range("F175").Value = 26
$ws.Range("G175").Value = 1807
$ws.Range("F179").Value = 55
$ws.Range("G179").Value = 3642.65
$ws.Range("B184").Value = 31714.85
$ws.Range("F199").Value = 49
$ws.Range("G199").Value = 1692.95
$ws.Range("F201").Value = 36
$ws.Range("G201").Value = 3672.36
$ws.Range("F202").Value = 178
$ws.Range("G202").Value = 12766.16
$ws.Range("F203").Value = 52
$ws.Range("G203").Value = 4632.68
$ws.Range("F206").Value = 17
$ws.Range("G206").Value = 801.55
$ws.Range("B208").Value = 34953.94
$ws.Range("F210").Value = 174
$ws.Range("G210").Value = 20344.08
$ws.Range("F211").Value = 1413
$ws.Range("G211").Value = 26140.5
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("F217").Value = 13
$ws.Range("G217").Value = 1815.32
$ws.Range("F221").Value = 25
$ws.Range("G221").Value = 1114.5
$ws.Range("B222").Value = 54751.05
$ws.Range("F237").Value = 8
$ws.Range("G237").Value = 256.56
$ws.Range("F239").Value = 2
$ws.Range("G239").Value = 164.06
$ws.Range("F241").Value = 63
$ws.Range("G241").Value = 18632.25
$ws.Range("F242").Value = 1
$ws.Range("G242").Value = 84.86
$ws.Range("F245").Value = 17
$ws.Range("G245").Value = 827.73
$ws.Range("F246").Value = 6
$ws.Range("G246").Value = 372.36
$ws.Range("F248").Value = 85
$ws.Range("G248").Value = 2725.95
$ws.Range("F250").Value = 6
$ws.Range("G250").Value = 281.22
$ws.Range("F255").Value = 15
$ws.Range("G255").Value = 8325.450000000001
$ws.Range("F256").Value = 0
$ws.Range("G256").Value = 0
$ws.Range("B258").Value = 44783.35
$ws.Range("F276").Value = 13
$ws.Range("G276").Value = 1663.74
$ws.Range("F284").Value = 78
$ws.Range("G284").Value = 11710.92
$ws.Range("F289").Value = 95
$ws.Range("G289").Value = 2509.9
$ws.Range("B290").Value = 70623.17
$ws.Range("F293").Value = 13
$ws.Range("G293").Value = 726.1799999999999
$ws.Range("B295").Value = 898.02
$ws.Range("F301").Value = 29
$ws.Range("G301").Value = 3256.12
$ws.Range("F304").Value = 0
$ws.Range("G304").Value = 0
$ws.Range("B307").Value = 7167.41
$ws.Range("F312").Value = 116
$ws.Range("G312").Value = 7012.2
$ws.Range("F315").Value = 7
$ws.Range("G315").Value = 292.18
$ws.Range("F329").Value = 176
$ws.Range("G329").Value = 16213.12
$ws.Range("F332").Value = 48
$ws.Range("G332").Value = 1934.4
$ws.Range("F333").Value = 11
$ws.Range("G333").Value = 174.35
$ws.Range("F337").Value = 114
$ws.Range("G337").Value = 8204.58
$ws.Range("F340").Value = 106
$ws.Range("G340").Value = 3050.68
$ws.Range("F342").Value = 49
$ws.Range("G342").Value = 3610.81
$ws.Range("F343").Value = 49
$ws.Range("G343").Value = 3103.17
$ws.Range("B347").Value = 122636.31
$ws.Range("F378").Value = 6
$ws.Range("G378").Value = 179.58
$ws.Range("F379").Value = 53
$ws.Range("G379").Value = 2832.32
$ws.Range("F387").Value = 73
$ws.Range("G387").Value = 1607.46
$ws.Range("F391").Value = 9
$ws.Range("G391").Value = 338.85
$ws.Range("B393").Value = 9016.200000000001
$ws.Range("F435").Value = 29
$ws.Range("G435").Value = 1481.32
$ws.Range("F437").Value = 6
$ws.Range("G437").Value = 727.02
$ws.Range("F439").Value = 40
$ws.Range("G439").Value = 2195.2
$ws.Range("F440").Value = 140
$ws.Range("G440").Value = 3910.2
$ws.Range("F442").Value = 97
$ws.Range("G442").Value = 2657.8
$ws.Range("B455").Value = 43622.06
$ws.Range("F537").Value = 0
$ws.Range("G537").Value = 0
$ws.Range("F541").Value = 27
$ws.Range("G541").Value = 2208.6
$ws.Range("F544").Value = 0
$ws.Range("G544").Value = 0
$ws.Range("F545").Value = 0
$ws.Range("G545").Value = 0
$ws.Range("F553").Value = 77
$ws.Range("G553").Value = 2641.87
$ws.Range("F554").Value = 0
$ws.Range("G554").Value = 0
$ws.Range("B555").Value = 7334.15
$ws.Range("F576").Value = 2
$ws.Range("G576").Value = 7226
$ws.Range("B578").Value = 8466.200000000001
$ws.Range("F596").Value = 75
$ws.Range("G596").Value = 3009.75
$ws.Range("F597").Value = 0
$ws.Range("G597").Value = 0
$ws.Range("F599").Value = 102
$ws.Range("G599").Value = 3970.86
$ws.Range("B601").Value = 20714.9
$ws.Range("B607").Value = 1587651.11
$ws.Range("B608").Value = 1587651.11

Write-Host "Applied 196 changes"